# feat: Enable AgGrid column grouping, and add import/export in customToolbar
#
# Data fixture update for cypress/fixtures/customers.xlsx:
#  - Howard's birthday corrected: 12/05/1987 -> 21/05/2002
#  - New row added for a customer named Billy (age 23, of legal drinking age,
#    favorite drink Beer, birthday 28 Apr 1940 entered as a real DATE()
#    formula, height 1.25)
#  - customers sheet becomes/stays the active sheet & selection moves to E7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers")

# --- fix Howard's birthday (row 5) ---
$ws.Range("E5").Value = "21/05/2002"

# --- append Billy as row 6 ---
$ws.Range("A6").Value = "Billy"
$ws.Range("B6").Value = 23

$ws.Range("C6").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("C6").Value = $true

$ws.Range("D6").Value = "Beer"

# Set the number format before writing the formula so Excel doesn't
# auto-register its own default "short date" format on top of ours.
$ws.Range("E6").NumberFormat = "mm/dd/yy"
$ws.Range("E6").Formula = "=DATE(1940,4,28)"

$ws.Range("F6").Value = 1.25

# --- selection / active sheet bookkeeping ---
$ws.Range("E7").Select() | Out-Null
